$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new top row (shifts header/items/total down by one row) ---
$ws.Rows.Item(1).Insert()

# --- New row 1: year headers for the cost columns ---
$ws.Range("B1").Value = 2017
$ws.Range("C1").Value = 2017
$ws.Range("D1").Value = 2021

# --- Row 2 (former row 1): extend the bold "Cost (Rands per month)" header
#     across the two new columns C and D ---
$ws.Range("C2").Value = "Cost (Rands per month)"
$ws.Range("D2").Value = "Cost (Rands per month)"
$ws.Range("C2").Font.Bold = $true
$ws.Range("D2").Font.Bold = $true

# --- Item rows 3-7 (former rows 2-6): fill the new per-year cost columns ---
$ws.Range("C3").Value = 320
$ws.Range("D3").Value = 419

$ws.Range("C4").Value = 950
$ws.Range("D4").Value = 1195

$ws.Range("C5").Value = 174.15
$ws.Range("D5").Value = 174.15

$ws.Range("C6").Value = 0
$ws.Range("D6").Value = 0

$ws.Range("C7").Value = 0
$ws.Range("D7").Value = 0

# --- Total row 9 (former row 8): totals for the new columns ---
$ws.Range("C9").Formula = "=SUM(C3:C8)"
$ws.Range("D9").Formula = "=SUM(D3:D8)"
$ws.Range("C9").Font.Bold = $true
$ws.Range("D9").Font.Bold = $true

# --- Match column widths of the new columns to column B ---
$ws.Range("C1").ColumnWidth = 19.33
$ws.Range("D1").ColumnWidth = 19.33

# --- Match the final selection state ---
[void]$ws.Range("D9").Select()
